$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Weeks of Work Required" for the engineering rows based on new prototype estimates
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 10
$ws.Range("C6").Value = 5

# Rename "Web Designer for Android" to "Web Designer for Web Application"
$ws.Range("A6").Value = "Web Designer for Web Application"
